$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "description"
$ws.Range("C2").Value = "Define your own project path where the code，data and result will be located"
$ws.Range("C3").Value = "If it put to 0, the project will run on a small dataset. If it put to 1, it will run on norman data set."

$ws.Columns.Item(3).ColumnWidth = 79.71

$ws.Range("C2").Select()
